$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting rows 8:124 down to 9:125
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with data
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44616
$ws.Range("D8").NumberFormat = $ws.Range("D9").NumberFormat
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 100112044
$ws.Range("G8").Value = "Perejil"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 2400
$ws.Range("K8").Value = 2300
$ws.Range("L8").Value = 2500
$ws.Range("M8").Value = 2400
$ws.Range("N8").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O8").Value = "Provincia del Elquí"
$ws.Range("P8").Value = 1600
$ws.Range("Q8").Value = 1.5
$ws.Range("R8").Value = "Hortaliza"
